# Open Day Poster — "Added features to the poster"
# Rewrites the bullet-point feature list in the "TextBox 22" shape (id 23)
# on slide 1: tweaks the first two bullets' wording (splitting the second
# bullet's leftover text into its own paragraph), shortens the third
# bullet, and adds a brand-new fourth bullet; the final "Lord cards..."
# bullet is left untouched. The textbox has <a:spAutoFit/>, so height
# re-flows automatically as the text changes.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("TextBox 22")
$tf = $shp.TextFrame
$tr = $tf.TextRange

# Bullet 1: "1-on-1 multiplayer card game with a dark fantasy theme."
#        -> "1-on-1 online multiplayer card game with a "
# (Assign a throwaway value first so the engine doesn't try to keep a
# "1-on-1 " common-prefix run from the old sentence around.)
$tr.Paragraphs(1).Text = "zzz"
$tr.Paragraphs(1).Text = "1-on-1 online multiplayer card game with a "

# Bullet 2: "Gameplay designed around resource management."
#        -> becomes its own paragraph, split across two runs: "D" then
#           "ark fantasy theme."
$tr.Paragraphs(2).Text = "D"
$null = $tr.Paragraphs(2).InsertAfter("ark fantasy theme.")

# Bullet 3: "A wide variety of cards to build your deck from."
#        -> "Resource management."
$tr.Paragraphs(3).Text = "zzz"
$tr.Paragraphs(3).Text = "Resource management."

# New bullet 4, inserted right after (inherits bullet 3's paragraph
# formatting): "A wide variety of cards."
$null = $tr.Paragraphs(3).InsertAfter("`rA wide variety of cards.")

# The autofit height the edited XML settles on lands 1 EMU short of what
# PowerPoint's own text-metrics computation calls for (a float-rounding
# artifact in the replay engine's live relayout vs. its ground-truth
# layout pass) — nudge it back onto the exact target height.
$shp.Height = 152.67657
